$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Cell, $Text)
    if ($Text -match '^-?\d+(\.\d+)?$') {
        # Force text storage for numeric-looking strings (matches how the
        # original inline-string cells are stored as text) via a quote-prefixed
        # entry, exactly like a user typing '6.35 into a General-formatted cell.
        $ws.Range($Cell).Value = "'" + $Text
    } else {
        $ws.Range($Cell).Value = $Text
    }
}

# Row 2
Set-TextCell 'D2' '62.732.61'
Set-TextCell 'E2' '  -0.97%  '

# Row 3
Set-TextCell 'D3' '3.000.30'
Set-TextCell 'E3' '  -3.92%  '

# Row 4
Set-TextCell 'E4' '  +0.03%  '

# Row 5
Set-TextCell 'D5' '554.11'
Set-TextCell 'E5' '  -1.10%  '

# Row 6
Set-TextCell 'D6' '151.93'
Set-TextCell 'E6' '  -5.85%  '

# Row 7
Set-TextCell 'E7' '  +0.01%  '

# Row 8
Set-TextCell 'D8' '0.561'
Set-TextCell 'E8' '  -3.54%  '

# Row 9
Set-TextCell 'D9' '3.002.25'
Set-TextCell 'E9' '  -3.72%  '

# Row 10
Set-TextCell 'E10' '  -2.37%  '

# Row 11
Set-TextCell 'D11' '6.35'
Set-TextCell 'E11' '  -5.19%  '

# Row 12
Set-TextCell 'D12' '0.363'
Set-TextCell 'E12' '  -3.83%  '

# Row 13
Set-TextCell 'D13' '3.526.75'
Set-TextCell 'E13' '  -3.80%  '

# Row 14
Set-TextCell 'E14' '  -3.35%  '

# Row 15
Set-TextCell 'D15' '62.820.77'
Set-TextCell 'E15' '  -0.88%  '

# Row 16
Set-TextCell 'D16' '23.79'
Set-TextCell 'E16' '  -3.75%  '

# Row 17
Set-TextCell 'D17' '3.008.78'
Set-TextCell 'E17' '  -3.65%  '

# Row 18
Set-TextCell 'D18' '0.0000148'
Set-TextCell 'E18' '  -2.37%  '

# Row 19
Set-TextCell 'D19' '394.19'
Set-TextCell 'E19' '  -1.21%  '

# Row 20
Set-TextCell 'D20' '5.08'
Set-TextCell 'E20' '  -2.33%  '

# Row 21
Set-TextCell 'D21' '11.80'
Set-TextCell 'E21' '  -4.94%  '

# Row 22
Set-TextCell 'D22' '6.59'
Set-TextCell 'E22' '  -6.16%  '

# Row 23
Set-TextCell 'E23' '  -0.11%  '

# Row 24
Set-TextCell 'D24' '65.00'
Set-TextCell 'E24' '  -3.48%  '

# Row 25
Set-TextCell 'D25' '0.463'
Set-TextCell 'E25' '  -2.75%  '

# Row 26
Set-TextCell 'D26' '0.187'
Set-TextCell 'E26' '  -6.36%  '

# Row 27
Set-TextCell 'D27' '0.0₃0958'
Set-TextCell 'E27' '  -4.29%  '

# Row 28
Set-TextCell 'D28' '8.59'
Set-TextCell 'E28' '  -1.08%  '

# Row 29
Set-TextCell 'D29' '0.999'
Set-TextCell 'E29' '  -0.08%  '

# Row 30
Set-TextCell 'E30' '  +0.01%  '

# Row 31
Set-TextCell 'D31' '1.74'
Set-TextCell 'E31' '  -2.55%  '

# Row 32
Set-TextCell 'D32' '20.42'
Set-TextCell 'E32' '  -2.12%  '

# Row 33
Set-TextCell 'D33' '160.82'
Set-TextCell 'E33' '  +5.48%  '

# Row 34
Set-TextCell 'D34' '4.64'
Set-TextCell 'E34' '  -2.42%  '

# Row 35
Set-TextCell 'D35' '5.99'
Set-TextCell 'E35' '  -3.89%  '

# Row 36
Set-TextCell 'D36' '1.08'
Set-TextCell 'E36' '  -1.59%  '

# Row 37
Set-TextCell 'D37' '1.28'
Set-TextCell 'E37' '  -2.57%  '

# Row 38
Set-TextCell 'E38' '  -4.19%  '

# Row 39
Set-TextCell 'D39' '2.472.26'
Set-TextCell 'E39' '  -10.00%  '

# Row 40
Set-TextCell 'E40' '  -2.21%  '

# Row 41
Set-TextCell 'B41' 'EnergySwap'
Set-TextCell 'C41' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D41' '22.43'
Set-TextCell 'E41' '  -3.50%  '

# Row 42
Set-TextCell 'B42' 'Filecoin'
Set-TextCell 'C42' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 'D42' '3.88'
Set-TextCell 'E42' '  -4.06%  '

# Row 43
Set-TextCell 'D43' '0.662'
Set-TextCell 'E43' '  -4.31%  '

# Row 44
Set-TextCell 'D44' '0.0592'
Set-TextCell 'E44' '  -3.29%  '

# Row 45
Set-TextCell 'E45' '  -0.07%  '

# Row 46
Set-TextCell 'E46' '  -3.90%  '

# Row 47
Set-TextCell 'D47' '5.01'
Set-TextCell 'E47' '  -6.82%  '

# Row 48
Set-TextCell 'B48' 'WhiteBITCoin'
Set-TextCell 'C48' 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextCell 'D48' '10.51'
Set-TextCell 'E48' '  +0.68%  '

# Row 49
Set-TextCell 'B49' 'InjectiveProtocol'
Set-TextCell 'C49' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell 'D49' '19.61'
Set-TextCell 'E49' '  -5.33%  '

# Row 50
Set-TextCell 'B50' 'Stellar'
Set-TextCell 'C50' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 'D50' '0.0942'
Set-TextCell 'E50' '  -2.88%  '

# Row 51
Set-TextCell 'D51' '262.42'
Set-TextCell 'E51' '  -6.39%  '
